$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (tab name) to reflect new "through" date
$ws.Name = "Through 2022-12-08"

# Update header label in column A, row 13 ("December (through 12-07)" -> "December (through 12-08)")
$ws.Range("A13").Value = "December (through 12-08)"

# Update I12 (2022 column, November row)
$ws.Range("I12").Value = 118

# Update row 13 (December) values for all year columns B..I
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 23
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 39
$ws.Range("H13").Value = 65
$ws.Range("I13").Value = 34

# Update row 14 (Total) values for all year columns B..I
$ws.Range("B14").Value = 299
$ws.Range("C14").Value = 586
$ws.Range("D14").Value = 851
$ws.Range("E14").Value = 699
$ws.Range("F14").Value = 544
$ws.Range("G14").Value = 1303
$ws.Range("H14").Value = 1708
$ws.Range("I14").Value = 1550
